# Update Convergence table for Tri M2: insert 6 new parameter rows
# (level_2[beta_lambda[0]], level_2[beta_mu[0]], level_2[beta_eta[0]],
#  level_2[beta_lambda[1]], level_2[beta_mu[1]], level_2[beta_eta[1]])
# right before the existing "level_2[var_log_lambda]" row, pushing the
# existing variance/covariance rows down by 6 rows (old rows 5-10 -> 11-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tri_Convergence_M2")

# Insert 6 blank rows starting at row 5 (shifts old rows 5:10 -> 11:16).
$ws.Range("A5:A10").EntireRow.Insert()

# The inserted rows don't inherit the bordered/bold/centered label style
# used by column A ("s=1" in the sheet XML). Copy that formatting from an
# existing labeled cell (A2) onto the new label cells.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A5:A10").PasteSpecial(-4122)

$newRows = @(
  @("level_2[beta_lambda[0]]", -3.8665, 0.2175, -4.2798, -3.4672, 0.0175, 0.0071, 153.3986, 481.8686, 1.0288),
  @("level_2[beta_mu[0]]",      0.0993, 0.1037, -0.0982,  0.2948, 0.0042, 0.0024, 614.7891, 1287.8694, 1.0081),
  @("level_2[beta_eta[0]]",     0.3791, 0.1938,  0.0135,  0.7495000000000001, 0.0086, 0.0058, 509.1998, 750.3375, 1.0048),
  @("level_2[beta_lambda[1]]",  3.2135, 0.08989999999999999, 3.0456, 3.3856, 0.0019, 0.0009, 2317.3609, 5136.2084, 1.0008),
  @("level_2[beta_mu[1]]",      0.0191, 0.1087, -0.1777,  0.2304, 0.002,  0.0023, 2903.5222, 4130.3821, 1.0061),
  @("level_2[beta_eta[1]]",     0.0361, 0.0315, -0.0242,  0.0936, 0.0005999999999999999, 0.0003, 2831.7283, 6337.3762, 1.0014)
)

$r = 5
foreach ($rowVals in $newRows) {
  for ($col = 1; $col -le $rowVals.Length; $col++) {
    $ws.Cells.Item($r, $col).Value = $rowVals[$col - 1]
  }
  $r++
}
